$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = "61.645.91"; Numeric = $false },
    @{ Cell = "E2"; Value = "  +1.33%  "; Numeric = $false },
    @{ Cell = "D3"; Value = "3.450.96"; Numeric = $false },
    @{ Cell = "E3"; Value = "  +2.11%  "; Numeric = $false },
    @{ Cell = "D4"; Value = "0.999"; Numeric = $true },
    @{ Cell = "E4"; Value = "  -0.12%  "; Numeric = $false },
    @{ Cell = "D5"; Value = "580.44"; Numeric = $true },
    @{ Cell = "E5"; Value = "  +1.49%  "; Numeric = $false },
    @{ Cell = "D6"; Value = "149.80"; Numeric = $true },
    @{ Cell = "E6"; Value = "  +9.13%  "; Numeric = $false },
    @{ Cell = "D7"; Value = "3.452.17"; Numeric = $false },
    @{ Cell = "E7"; Value = "  +2.22%  "; Numeric = $false },
    @{ Cell = "E8"; Value = "  +0.04%  "; Numeric = $false },
    @{ Cell = "E9"; Value = "  +0.97%  "; Numeric = $false },
    @{ Cell = "E10"; Value = "  +1.79%  "; Numeric = $false },
    @{ Cell = "E11"; Value = "  +2.87%  "; Numeric = $false },
    @{ Cell = "D13"; Value = "4.036.39"; Numeric = $false },
    @{ Cell = "E13"; Value = "  +1.95%  "; Numeric = $false },
    @{ Cell = "D14"; Value = "28.02"; Numeric = $true },
    @{ Cell = "E14"; Value = "  +7.80%  "; Numeric = $false },
    @{ Cell = "E15"; Value = "  -0.29%  "; Numeric = $false },
    @{ Cell = "E16"; Value = "  +1.72%  "; Numeric = $false },
    @{ Cell = "D17"; Value = "3.449.43"; Numeric = $false },
    @{ Cell = "E17"; Value = "  +1.70%  "; Numeric = $false },
    @{ Cell = "D18"; Value = "61.754.22"; Numeric = $false },
    @{ Cell = "E18"; Value = "  +1.17%  "; Numeric = $false },
    @{ Cell = "D19"; Value = "6.29"; Numeric = $true },
    @{ Cell = "E19"; Value = "  +8.46%  "; Numeric = $false },
    @{ Cell = "E20"; Value = "  +2.94%  "; Numeric = $false },
    @{ Cell = "D21"; Value = "9.52"; Numeric = $true },
    @{ Cell = "E21"; Value = "  +0.84%  "; Numeric = $false },
    @{ Cell = "D22"; Value = "390.28"; Numeric = $true },
    @{ Cell = "E22"; Value = "  +4.17%  "; Numeric = $false },
    @{ Cell = "E23"; Value = "  +2.39%  "; Numeric = $false },
    @{ Cell = "D24"; Value = "3.588.37"; Numeric = $false },
    @{ Cell = "E24"; Value = "  +1.72%  "; Numeric = $false },
    @{ Cell = "D25"; Value = "73.00"; Numeric = $true },
    @{ Cell = "E25"; Value = "  +2.74%  "; Numeric = $false },
    @{ Cell = "E26"; Value = "  +0.09%  "; Numeric = $false },
    @{ Cell = "E27"; Value = "  +0.68%  "; Numeric = $false },
    @{ Cell = "E28"; Value = "  +0.29%  "; Numeric = $false },
    @{ Cell = "E29"; Value = "  +4.20%  "; Numeric = $false },
    @{ Cell = "E30"; Value = "  +3.56%  "; Numeric = $false },
    @{ Cell = "E31"; Value = "  -12.11%  "; Numeric = $false },
    @{ Cell = "E32"; Value = "  +0.89%  "; Numeric = $false },
    @{ Cell = "D33"; Value = "8.27"; Numeric = $true },
    @{ Cell = "E33"; Value = "  +1.63%  "; Numeric = $false },
    @{ Cell = "E34"; Value = "  +1.48%  "; Numeric = $false },
    @{ Cell = "D36"; Value = "24.09"; Numeric = $true },
    @{ Cell = "E36"; Value = "  +1.89%  "; Numeric = $false },
    @{ Cell = "E37"; Value = "  +2.84%  "; Numeric = $false },
    @{ Cell = "D38"; Value = "5.21"; Numeric = $true },
    @{ Cell = "E38"; Value = "  +0.38%  "; Numeric = $false },
    @{ Cell = "E39"; Value = "  +1.31%  "; Numeric = $false },
    @{ Cell = "D40"; Value = "166.96"; Numeric = $true },
    @{ Cell = "D41"; Value = "0.0789"; Numeric = $true },
    @{ Cell = "E41"; Value = "  +3.91%  "; Numeric = $false },
    @{ Cell = "D42"; Value = "27.15"; Numeric = $true },
    @{ Cell = "E42"; Value = "  +12.59%  "; Numeric = $false },
    @{ Cell = "D43"; Value = "0.795"; Numeric = $true },
    @{ Cell = "E43"; Value = "  +2.36%  "; Numeric = $false },
    @{ Cell = "E44"; Value = "  +2.33%  "; Numeric = $false },
    @{ Cell = "D45"; Value = "0.999"; Numeric = $true },
    @{ Cell = "E45"; Value = "  -0.16%  "; Numeric = $false },
    @{ Cell = "D46"; Value = "42.38"; Numeric = $true },
    @{ Cell = "E46"; Value = "  +1.83%  "; Numeric = $false },
    @{ Cell = "E47"; Value = "  +0.51%  "; Numeric = $false },
    @{ Cell = "D48"; Value = "2.601.15"; Numeric = $false },
    @{ Cell = "E48"; Value = "  +6.10%  "; Numeric = $false },
    @{ Cell = "D49"; Value = "1.16"; Numeric = $true },
    @{ Cell = "E49"; Value = "  -2.79%  "; Numeric = $false },
    @{ Cell = "D50"; Value = "6.97"; Numeric = $true },
    @{ Cell = "E50"; Value = "  +2.50%  "; Numeric = $false },
    @{ Cell = "D51"; Value = "23.27"; Numeric = $true },
    @{ Cell = "E51"; Value = "  +1.04%  "; Numeric = $false }
)

foreach ($c in $changes) {
    if ($c.Numeric) {
        $ws.Range($c.Cell).NumberFormat = "@"
        $ws.Range($c.Cell).Value = $c.Value
        $ws.Range($c.Cell).Style = "Normal"
    } else {
        $ws.Range($c.Cell).Value = $c.Value
    }
}

Write-Host "Applied $($changes.Count) cell updates"
